# Swap the taxon-observation data between rows 13 and 14 in the "Artfynd" sheet.
# Columns affected: A, B, E, F, G, H, P, Q, R, S
# (All other columns, e.g. C, D, I, T..AY stay attached to their original row.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")

foreach ($col in $cols) {
    $cell13 = $ws.Range($col + "13")
    $cell14 = $ws.Range($col + "14")

    $val13 = $cell13.Value()
    $val14 = $cell14.Value()

    $cell13.Value = $val14
    $cell14.Value = $val13
}
